$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number by Excel's type
# inference (e.g. "1.00", "28.16") get their NumberFormat pinned to
# "@" (text) before the write, then their original Style restored so
# the saved file keeps the original (default) style index.
function Set-TextValue {
    param($cellRef, $text)
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '60.258.25'
Set-TextValue 'D3' '3.378.16'
$ws.Range('E3').Value = '  -2.09%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '571.38'
$ws.Range('E5').Value = '  -1.43%  '
Set-TextValue 'D6' '141.47'
$ws.Range('E6').Value = '  -4.56%  '
$ws.Range('E7').Value = '  +0.04%  '
Set-TextValue 'D8' '3.380.21'
$ws.Range('E8').Value = '  -2.05%  '
Set-TextValue 'D9' '0.475'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('E10').Value = '  -4.13%  '
$ws.Range('E11').Value = '  -0.63%  '
Set-TextValue 'D12' '0.395'
$ws.Range('E12').Value = '  +0.80%  '
Set-TextValue 'D13' '3.959.49'
$ws.Range('E13').Value = '  -2.05%  '
Set-TextValue 'D14' '28.16'
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('E16').Value = '  -2.26%  '
Set-TextValue 'D17' '3.384.00'
$ws.Range('E17').Value = '  -2.11%  '
Set-TextValue 'D18' '60.435.21'
$ws.Range('E18').Value = '  -1.93%  '
Set-TextValue 'D19' '6.27'
$ws.Range('E19').Value = '  -0.89%  '
Set-TextValue 'D20' '14.09'
$ws.Range('E20').Value = '  -1.77%  '
Set-TextValue 'D21' '9.12'
$ws.Range('E21').Value = '  -3.52%  '
Set-TextValue 'D22' '389.07'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('E23').Value = '  -1.74%  '
Set-TextValue 'D24' '73.36'
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('E26').Value = '  -4.46%  '
Set-TextValue 'D27' '3.519.40'
$ws.Range('E27').Value = '  -2.18%  '
Set-TextValue 'D28' '0.178'
$ws.Range('E28').Value = '  -1.03%  '
Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  -5.50%  '
Set-TextValue 'D31' '8.07'
$ws.Range('E31').Value = '  -2.25%  '
Set-TextValue 'D32' '2.14'
$ws.Range('E32').Value = '  -1.28%  '
Set-TextValue 'D33' '1.43'
$ws.Range('E33').Value = '  -6.55%  '
$ws.Range('E34').Value = '  -0.06%  '
Set-TextValue 'D35' '23.74'
$ws.Range('E35').Value = '  -0.99%  '
Set-TextValue 'D36' '6.94'
$ws.Range('E36').Value = '  -1.79%  '
Set-TextValue 'D37' '3.411.09'
$ws.Range('E37').Value = '  -1.83%  '
Set-TextValue 'D38' '167.07'
$ws.Range('E38').Value = '  +0.36%  '
Set-TextValue 'D39' '4.98'
$ws.Range('E39').Value = '  -4.97%  '
$ws.Range('E40').Value = '  -3.71%  '
Set-TextValue 'D41' '0.0778'
$ws.Range('E41').Value = '  -1.58%  '
Set-TextValue 'D42' '27.04'
$ws.Range('E42').Value = '  +3.29%  '
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D46' '1.69'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D47' '41.12'
$ws.Range('E47').Value = '  -2.81%  '
Set-TextValue 'D48' '2.529.91'
$ws.Range('E48').Value = '  -2.96%  '
$ws.Range('E49').Value = '  -3.08%  '
Set-TextValue 'D50' '6.85'
$ws.Range('E50').Value = '  -1.82%  '
Set-TextValue 'D51' '23.03'
$ws.Range('E51').Value = '  -1.10%  '
